# Renamed a few transcripts / updated the DataSheet:
#   - Column D (Speaker): "RT1" -> "T1", "Students" -> "SS", "Student" -> "S"
#   - Column F (Teacher Tag): "3 - getting students to relate" -> "3 - getting SS to relate"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $speakerCell = $ws.Cells.Item($r, 4)
    $speakerVal = $speakerCell.Value2

    if ($speakerVal -eq "RT1") {
        $speakerCell.Value = "T1"
    } elseif ($speakerVal -eq "Students") {
        $speakerCell.Value = "SS"
    } elseif ($speakerVal -eq "Student") {
        $speakerCell.Value = "S"
    }

    $tagCell = $ws.Cells.Item($r, 6)
    $tagVal = $tagCell.Value2

    if ($tagVal -eq "3 - getting students to relate") {
        $tagCell.Value = "3 - getting SS to relate"
    }
}
